$d = $word.ActiveDocument

# Original:
#   "...webapplikation 'Budgetmanager' som er en app virksomheder kan bruge..."
# New:
#   "...webapplikation 'Budgetmanager' der er en app som virksomheder kan bruge..."
#
# Use a fairly long, unique anchor of surrounding text so the replace is
# unambiguous and only touches this one spot in the document.
$d.Content.Find.Execute(
    "Budgetmanager’ som er en app virksomheder",
    $true,                                    # MatchCase
    $false,                                   # MatchWholeWord
    $false,                                   # MatchWildcards
    $false,                                   # MatchSoundsLike
    $false,                                   # MatchAllWordForms
    $true,                                    # Forward
    1,                                        # Wrap (wdFindContinue)
    $false,                                   # Format
    "Budgetmanager’ der er en app som virksomheder",
    2                                         # Replace (wdReplaceAll)
) | Out-Null
